# Update the AC cable existing capacity matrix (Sheet1, range B2:I9) with
# newly computed cluster-to-cluster capacity values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - cluster_0
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 7
$ws.Range("D2").Value = 11
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 18
$ws.Range("G2").Value = 4
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 2.5

# Row 3 - cluster_1
$ws.Range("B3").Value = 7
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 10
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 3
$ws.Range("H3").Value = 6
$ws.Range("I3").Value = 0

# Row 4 - cluster_2
$ws.Range("B4").Value = 11
$ws.Range("C4").Value = 10
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 3

# Row 5 - cluster_3
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 12
$ws.Range("H5").Value = 2.4
$ws.Range("I5").Value = 0

# Row 6 - cluster_4
$ws.Range("B6").Value = 18
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 4
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0

# Row 7 - cluster_5
$ws.Range("B7").Value = 4
$ws.Range("C7").Value = 3
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 12
$ws.Range("F7").Value = 4
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0

# Row 8 - cluster_6
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 6
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 2.4
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0

# Row 9 - cluster_7
$ws.Range("B9").Value = 2.5
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 3
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
